$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2172619047619048
$ws.Range("C2").Value = 0.4970238095238095
$ws.Range("J2").Value = 0.03571428571428571
$ws.Range("P2").Value = 0.1398809523809524
$ws.Range("S2").Value = 0.1101190476190476

# Row 3
$ws.Range("B3").Value = 0.01162790697674419
$ws.Range("C3").Value = 0.04069767441860465
$ws.Range("J3").Value = 0.04651162790697674
$ws.Range("P3").Value = 0.6104651162790697
$ws.Range("S3").Value = 0.2906976744186047

# Row 4
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("O4").Value = 0.0196078431372549
$ws.Range("P4").Value = 0.7450980392156863
$ws.Range("S4").Value = 0.1764705882352941

# Row 6
$ws.Range("B6").Value = 0.06936416184971098
$ws.Range("D6").Value = 0.01734104046242774
$ws.Range("F6").Value = 0.03468208092485549
$ws.Range("J6").Value = 0.2890173410404624
$ws.Range("O6").Value = 0.03468208092485549
$ws.Range("Q6").Value = 0.1734104046242775
$ws.Range("R6").Value = 0.06936416184971098
$ws.Range("S6").Value = 0.3121387283236994

# Row 7
$ws.Range("B7").Value = 0.1345029239766082
$ws.Range("D7").Value = 0.02923976608187134
$ws.Range("E7").Value = 0.005847953216374269
$ws.Range("F7").Value = 0.04093567251461988
$ws.Range("J7").Value = 0.2105263157894737
$ws.Range("O7").Value = 0.01169590643274854
$ws.Range("Q7").Value = 0.1578947368421053
$ws.Range("R7").Value = 0.06432748538011696
$ws.Range("S7").Value = 0.3450292397660819

# Row 8
$ws.Range("B8").Value = 0.08542713567839195
$ws.Range("D8").Value = 0.02512562814070352
$ws.Range("E8").Value = 0.002512562814070352
$ws.Range("F8").Value = 0.04522613065326633
$ws.Range("J8").Value = 0.1909547738693467
$ws.Range("O8").Value = 0.01758793969849246
$ws.Range("Q8").Value = 0.1884422110552764
$ws.Range("R8").Value = 0.08291457286432161
$ws.Range("S8").Value = 0.3618090452261307

# Row 9
$ws.Range("B9").Value = 0.1271676300578035
$ws.Range("D9").Value = 0.005780346820809248
$ws.Range("F9").Value = 0.02890173410404624
$ws.Range("J9").Value = 0.138728323699422
$ws.Range("O9").Value = 0.02890173410404624
$ws.Range("Q9").Value = 0.2543352601156069
$ws.Range("R9").Value = 0.07514450867052024
$ws.Range("S9").Value = 0.3410404624277457

# Row 10
$ws.Range("B10").Value = 0.1216702663786897
$ws.Range("D10").Value = 0.02447804175665947
$ws.Range("F10").Value = 0.06047516198704104
$ws.Range("J10").Value = 0.1547876169906408
$ws.Range("O10").Value = 0.02375809935205184
$ws.Range("Q10").Value = 0.1943844492440605
$ws.Range("R10").Value = 0.08063354931605471
$ws.Range("S10").Value = 0.339812814974802

# Row 11
$ws.Range("G11").Value = 0.1486988847583643
$ws.Range("J11").Value = 0.1078066914498141
$ws.Range("K11").Value = 0.2007434944237918
$ws.Range("L11").Value = 0.5315985130111525
$ws.Range("S11").Value = 0.01115241635687732

# Row 12
$ws.Range("G12").Value = 0.6470588235294118
$ws.Range("J12").Value = 0.2745098039215687
$ws.Range("K12").Value = 0.0196078431372549
$ws.Range("L12").Value = 0.0457516339869281
$ws.Range("S12").Value = 0.0130718954248366

# Row 13
$ws.Range("G13").Value = 0.7555555555555555
$ws.Range("J13").Value = 0.2444444444444444

# Row 15
$ws.Range("F15").Value = 0.01731601731601732
$ws.Range("H15").Value = 0.1212121212121212
$ws.Range("I15").Value = 0.06493506493506493
$ws.Range("J15").Value = 0.3896103896103896
$ws.Range("K15").Value = 0.05627705627705628
$ws.Range("M15").Value = 0.004329004329004329
$ws.Range("O15").Value = 0.05627705627705628
$ws.Range("S15").Value = 0.29004329004329

# Row 16
$ws.Range("F16").Value = 0.02150537634408602
$ws.Range("H16").Value = 0.1774193548387097
$ws.Range("I16").Value = 0.06451612903225806
$ws.Range("J16").Value = 0.4462365591397849
$ws.Range("K16").Value = 0.1182795698924731
$ws.Range("M16").Value = 0.005376344086021506
$ws.Range("N16").Value = 0.005376344086021506
$ws.Range("O16").Value = 0.05376344086021505
$ws.Range("S16").Value = 0.1075268817204301

# Row 17
$ws.Range("H17").Value = 0.1805869074492099
$ws.Range("I17").Value = 0.07223476297968397
$ws.Range("J17").Value = 0.4604966139954853
$ws.Range("K17").Value = 0.0835214446952596
$ws.Range("M17").Value = 0.01580135440180587
$ws.Range("O17").Value = 0.06320541760722348
$ws.Range("S17").Value = 0.1241534988713318

# Row 18
$ws.Range("F18").Value = 0.00558659217877095
$ws.Range("H18").Value = 0.2067039106145251
$ws.Range("I18").Value = 0.0893854748603352
$ws.Range("J18").Value = 0.4301675977653631
$ws.Range("K18").Value = 0.09497206703910614
$ws.Range("M18").Value = 0.00558659217877095
$ws.Range("O18").Value = 0.0670391061452514
$ws.Range("S18").Value = 0.1005586592178771

# Row 19
$ws.Range("F19").Value = 0.01197604790419162
$ws.Range("H19").Value = 0.1873396065012831
$ws.Range("I19").Value = 0.08297690333618478
$ws.Range("J19").Value = 0.3875106928999145
$ws.Range("K19").Value = 0.1026518391787853
$ws.Range("M19").Value = 0.02822925577416596
$ws.Range("N19").Value = 0.002566295979469632
$ws.Range("O19").Value = 0.07869974337040206
$ws.Range("S19").Value = 0.1180496150556031
